$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-08-26 00:26:50", "geneRNASeq", "stress", 40, 8, 32, 0.04003624967153208),
    @("2023-08-26 00:26:50", "geneRNASeq", "stress", 10, 3, 7, 0.104927894989448),
    @("2023-08-26 00:26:51", "geneRNASeq", "stress", 20, 5, 15, 0.06366885883235399),
    @("2023-08-26 00:26:51", "geneRNASeq", "stress", 30, 8, 22, 0.04828553388125874)
)

$startRow = 25
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
    $ws.Cells.Item($row, 7).Value = $rowData[6]
}
